{"js": "// Replace the date line and each \"A\u00d7B=C\" answer cell with its new value,\n// per the commit's updated output (new random problem set + date bump).\nconst replacements = [\n  [\"2024-03-14 Thursday\", \"2024-03-15 Friday\"],\n  [\"817\u00d78=6536\", \"781\u00d73=2343\"],\n  [\"795\u00d79=7155\", \"911\u00d72=1822\"],\n  [\"243\u00d73=729\", \"959\u00d76=5754\"],\n  [\"832\u00d78=6656\", \"963\u00d74=3852\"],\n  [\"531\u00d75=2655\", \"206\u00d72=412\"],\n  [\"705\u00d74=2820\", \"166\u00d76=996\"],\n  [\"705\u00d79=6345\", \"230\u00d74=920\"],\n  [\"162\u00d75=810\", \"756\u00d76=4536\"],\n  [\"727\u00d72=1454\", \"966\u00d72=1932\"],\n  [\"983\u00d75=4915\", \"341\u00d72=682\"],\n  [\"716\u00d79=6444\", \"201\u00d77=1407\"],\n  [\"935\u00d78=7480\", \"803\u00d77=5621\"],\n  [\"969\u00d78=7752\", \"543\u00d79=4887\"],\n  [\"699\u00d75=3495\", \"275\u00d72=550\"],\n  [\"291\u00d75=1455\", \"287\u00d78=2296\"],\n  [\"232\u00d77=1624\", \"693\u00d78=5544\"],\n  [\"236\u00d74=944\", \"585\u00d74=2340\"],\n  [\"157\u00d73=471\", \"932\u00d76=5592\"],\n  [\"419\u00d79=3771\", \"154\u00d78=1232\"],\n  [\"746\u00d75=3730\", \"522\u00d72=1044\"],\n  [\"808\u00d75=4040\", \"900\u00d75=4500\"],\n  [\"200\u00d74=800\", \"478\u00d76=2868\"],\n  [\"863\u00d74=3452\", \"301\u00d78=2408\"],\n  [\"662\u00d72=1324\", \"382\u00d76=2292\"],\n  [\"562\u00d72=1124\", \"364\u00d74=1456\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"A\u00d7B=C\" answer cell with its new value,\n# per the commit's updated output (new random problem set + date bump).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-14 Thursday\", \"2024-03-15 Friday\"),\n    @(\"817\u00d78=6536\", \"781\u00d73=2343\"),\n    @(\"795\u00d79=7155\", \"911\u00d72=1822\"),\n    @(\"243\u00d73=729\", \"959\u00d76=5754\"),\n    @(\"832\u00d78=6656\", \"963\u00d74=3852\"),\n    @(\"531\u00d75=2655\", \"206\u00d72=412\"),\n    @(\"705\u00d74=2820\", \"166\u00d76=996\"),\n    @(\"705\u00d79=6345\", \"230\u00d74=920\"),\n    @(\"162\u00d75=810\", \"756\u00d76=4536\"),\n    @(\"727\u00d72=1454\", \"966\u00d72=1932\"),\n    @(\"983\u00d75=4915\", \"341\u00d72=682\"),\n    @(\"716\u00d79=6444\", \"201\u00d77=1407\"),\n    @(\"935\u00d78=7480\", \"803\u00d77=5621\"),\n    @(\"969\u00d78=7752\", \"543\u00d79=4887\"),\n    @(\"699\u00d75=3495\", \"275\u00d72=550\"),\n    @(\"291\u00d75=1455\", \"287\u00d78=2296\"),\n    @(\"232\u00d77=1624\", \"693\u00d78=5544\"),\n    @(\"236\u00d74=944\", \"585\u00d74=2340\"),\n    @(\"157\u00d73=471\", \"932\u00d76=5592\"),\n    @(\"419\u00d79=3771\", \"154\u00d78=1232\"),\n    @(\"746\u00d75=3730\", \"522\u00d72=1044\"),\n    @(\"808\u00d75=4040\", \"900\u00d75=4500\"),\n    @(\"200\u00d74=800\", \"478\u00d76=2868\"),\n    @(\"863\u00d74=3452\", \"301\u00d78=2408\"),\n    @(\"662\u00d72=1324\", \"382\u00d76=2292\"),\n    @(\"562\u00d72=1124\", \"364\u00d74=1456\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n"}
